$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row before row 3 ---
# This shifts the old rows 3..9 down to 4..10 and shifts all merges/styles
# that reference those rows automatically.
$ws.Rows(3).Insert()

# --- Turn the (now two-row-tall) row2/row3 label cells into merged cells ---
# A2/B2/C2 already hold "姓名*" / "联系方式*" / "收件地址*"; merge them down
# into the freshly inserted row 3 so each label spans both rows.
$ws.Range("A2:A3").Merge()
$ws.Range("B2:B3").Merge()
$ws.Range("C2:C3").Merge()

# Copy the row2 formatting (font/border/alignment) down onto row3 so the
# newly inserted row looks consistent with the label row above it.
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Add the new "手机号码*：" / "QQ*：" labels in column D ---
$ws.Range("C2").Copy()
$ws.Range("D2:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D2:D3").HorizontalAlignment = -4131
$ws.Range("D2").Value = "手机号码*："
$ws.Range("D3").Value = "QQ*："

# --- Row heights ---
# (rows 5 & 6 are left untouched - their height was never a custom override
# either before or after the edit, it is just the natural 20pt row height)
$ws.Rows(1).RowHeight = 25
$ws.Rows(2).RowHeight = 25
$ws.Rows(3).RowHeight = 25
$ws.Rows(4).RowHeight = 45
$ws.Rows(7).RowHeight = 40
$ws.Rows(8).RowHeight = 35
$ws.Rows(9).RowHeight = 45
$ws.Rows(10).RowHeight = 15

# --- Column widths ---
$ws.Columns(1).ColumnWidth = 19.5
$ws.Columns(2).ColumnWidth = 14.5
$ws.Columns(3).ColumnWidth = 19.5
$ws.Columns(4).ColumnWidth = 29.5

# --- View: zoom + selection ---
$ws.Range("A4").Select()
$excel.ActiveWindow.Zoom = 62

Write-Host "edit complete"
